$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.298.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "1.586.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "  +1.15%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'213.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'23.96"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.75%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.250"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0886"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "1.812.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "1.597.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "28.278.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'63.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'227.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "0.0₃0705"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'151.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'15.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "1.397.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.06%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "  -7.89%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "  +8.99%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.539"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.810"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.979"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'64.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "1.721.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'86.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  +10.55%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "  -0.81%  "
$ws.Range("E51").Style = "Normal"
